$wb = $excel.ActiveWorkbook

# --- Add the two new sheets at the end of the workbook -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLifetimeCategories = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsLifetimeCategories.Name = "lifetime_categories"

$wsLifetimeBreakdowns = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLifetimeCategories)
$wsLifetimeBreakdowns.Name = "lifetime_breakdowns"

# --- Populate lifetime_categories -----------------------------------------
# Order matters: it reproduces the shared-string table order Excel would
# produce from the authoring sequence (new strings appended in first-use
# order, old unreferenced "All Other" pruned on save).
$wsLifetimeCategories.Range("A1").Value = "First house downpayment"

# --- Rename "All Other" -> "All other" on the categories sheet ------------
$wsCategories = $wb.Worksheets.Item("categories")
$wsCategories.Range("A8").Value = "All other"

$wsLifetimeCategories.Range("A2").Value = "Pre-school childcare"
$wsLifetimeCategories.Range("B1").Value = "House"
$wsLifetimeCategories.Range("B2").Value = "Childcare"
$wsLifetimeCategories.Range("C1").Value = "Deposit for buying a first house"
$wsLifetimeCategories.Range("C2").Value = "Monthly payments for childcare, over and above government-provided childcare hours of 15-30 hours/week during term times for ages 3 and over"
$wsLifetimeCategories.Range("B3").Value = "Childcareyears"
$wsLifetimeCategories.Range("A3").Value = "Pre-school childcare"
$wsLifetimeCategories.Range("C3").Value = "Monthly payments for childcare, over and above government-provided childcare hours of 15-30 hours/week during term times for ages 3 and over"

# --- Formatting for lifetime_categories ------------------------------------
# (Best-fit column width on column A for the longest label it holds.)
$wsLifetimeCategories.Columns.Item(1).ColumnWidth = 22.92

# --- Populate lifetime_breakdowns ------------------------------------------
$wsLifetimeBreakdowns.Range("A1").Value = 20000
$wsLifetimeBreakdowns.Range("B1").Value = 600
$wsLifetimeBreakdowns.Range("C1").Value = 2

$wsLifetimeBreakdowns.Range("A2").Value = 20000
$wsLifetimeBreakdowns.Range("B2").Value = 600
$wsLifetimeBreakdowns.Range("C2").Value = 2

$wsLifetimeBreakdowns.Range("A3").Value = 20000
$wsLifetimeBreakdowns.Range("B3").Value = 600
$wsLifetimeBreakdowns.Range("C3").Value = 2

$wsLifetimeBreakdowns.Range("A4").Value = 20000
$wsLifetimeBreakdowns.Range("B4").Value = 600
$wsLifetimeBreakdowns.Range("C4").Value = 2

$wsLifetimeBreakdowns.Range("A5").Value = 20000
$wsLifetimeBreakdowns.Range("B5").Value = 600
$wsLifetimeBreakdowns.Range("C5").Value = 2

$wsLifetimeBreakdowns.Range("A6").Value = 50000
$wsLifetimeBreakdowns.Range("B6").Value = 1200
$wsLifetimeBreakdowns.Range("C6").Value = 4

# --- Selections -------------------------------------------------------------
$wsCategories.Activate()
$wsCategories.Range("A9").Select() | Out-Null

$wsLifetimeBreakdowns.Activate()
$wsLifetimeBreakdowns.Range("D10").Select() | Out-Null

$wsLifetimeCategories.Activate()
$wsLifetimeCategories.Range("C2").Select() | Out-Null
